# Update town close year columns
# Rename the year-specific header labels in row 1 of Sheet1 from hard-coded
# years (2023/2024) to relative "Prior Year" / "Curr. Year" labels so the
# template doesn't need to be edited every tax year.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("O1").Value = "Curr. Year LMV"
$ws.Range("P1").Value = "Curr. Year BMV"
$ws.Range("Q1").Value = "Curr. Year TMV"
$ws.Range("R1").Value = "Prior Year LMV"
$ws.Range("S1").Value = "Prior Year BMV"
$ws.Range("T1").Value = "Prior Year TMV"
$ws.Range("V1").Value = "Curr. Year Dwelling MV"
$ws.Range("W1").Value = "Curr. Year Dwelling Total"
$ws.Range("X1").Value = "Prior Year Dwelling MV"
$ws.Range("Y1").Value = "Prior Year Dwelling Total"
